$d = $word.ActiveDocument

function New-WordPackageXml($innerBodyXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerBodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph near the top of the document.
#    It immediately follows the "Play Christmas Gift Rush for Free - Review"
#    heading, and contains a bold "Meta description" run plus the
#    description-text run.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Meta description*Get into the holiday spirit*") {
        [void]$para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Insert a brand-new paragraph ("Play Christmas Gift Rush for Free -
#    Review", bold) right before the final "Prompt: Create a cartoon..."
#    paragraph.
# ---------------------------------------------------------------------------
$promptIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Prompt: Create a cartoon style feature image*") {
        $promptIdx = $i
        break
    }
}

if ($promptIdx -gt 0) {
    $promptPara = $d.Paragraphs.Item($promptIdx)
    $insertionPoint = $promptPara.Range.Duplicate
    $insertionPoint.Collapse(1)

    $newParaXml = New-WordPackageXml('<w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Christmas Gift Rush for Free - Review</w:t></w:r></w:p></w:body>')
    [void]$insertionPoint.InsertXML($newParaXml)

    # The prompt paragraph index shifted down by one because of the new
    # paragraph inserted just before it.
    $promptIdx = $promptIdx + 1
}

# ---------------------------------------------------------------------------
# 3. Replace the text of the (now shifted) "Prompt: ..." paragraph with the
#    meta-description text, keeping its existing italic run formatting and
#    leading empty run.
# ---------------------------------------------------------------------------
if ($promptIdx -gt 0) {
    $promptPara = $d.Paragraphs.Item($promptIdx)

    $textRange = $promptPara.Range.Duplicate
    $textRange.MoveEnd(1, -1)
    $textRange.Delete()

    $insertionPoint = $d.Paragraphs.Item($promptIdx).Range.Duplicate
    $insertionPoint.Collapse(1)

    $newTextXml = New-WordPackageXml('<w:body><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Get into the holiday spirit with Christmas Gift Rush! Read our review and play for free. Enjoy excellent graphics, high payouts, and a unique Nudge function.</w:t></w:r></w:p></w:body>')
    $insertionPoint.InsertXML($newTextXml)

    # InsertXML inserted a full paragraph (with its own paragraph mark),
    # leaving the paragraph mark that used to belong to the old "Prompt"
    # text as a stray trailing empty paragraph. Merge it back by deleting
    # the paragraph mark that separates the two.
    $mergedPara = $d.Paragraphs.Item($promptIdx)
    $markRange = $mergedPara.Range.Duplicate
    $markRange.Collapse(0)
    $markRange.MoveEnd(1, 1)
    if ($markRange.Text.Length -gt 0) {
        $markRange.Delete()
    }
}
